$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Return(percent) statistics)
$ws.Range("B2").Value = -0.0001234629060598564
$ws.Range("D2").Value = -0.0995249943451703
$ws.Range("G2").Value = 0.013131976362443
$ws.Range("H2").Value = 0.09784051995248119

# Row 3 (Volume(1,000 shares) statistics)
$ws.Range("B3").Value = 411.0832039899008
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 41.75012
$ws.Range("E3").Value = 149.017
$ws.Range("F3").Value = 220.378
$ws.Range("G3").Value = 334.861
$ws.Range("H3").Value = 1938.810599999995
